$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated experiment results for random_forest row (new k for cross-validation)
$ws.Range("B2").Value = 53.857008191017108
$ws.Range("C2").Value = 0.53551799274112422
$ws.Range("D2").Value = 38.300452914555407
$ws.Range("E2").Value = 0.40660222850719541
$ws.Range("F2").Value = 0.63765369010709516
$ws.Range("G2").Value = 0.58895856792883605
$ws.Range("H2").Value = 0.59339777149280459
$ws.Range("I2").Value = 0.77094371092881409

# Updated experiment results for lsboost row (new k for cross-validation)
$ws.Range("B3").Value = 52.876684427421708
$ws.Range("C3").Value = 0.52577031028065291
$ws.Range("D3").Value = 35.909416643758128
$ws.Range("E3").Value = 0.39193471924465073
$ws.Range("F3").Value = 0.62604689859838036
$ws.Range("G3").Value = 0.5521908226215867
$ws.Range("H3").Value = 0.60806528075534927
$ws.Range("I3").Value = 0.78305083475347004

# Updated experiment results for neural_network row (new k for cross-validation)
$ws.Range("B4").Value = 57.191976264880665
$ws.Range("C4").Value = 0.56867868006406064
$ws.Range("D4").Value = 38.288827173006716
$ws.Range("E4").Value = 0.45851707546108361
$ws.Range("F4").Value = 0.67713888934330424
$ws.Range("G4").Value = 0.58877979510573397
$ws.Range("H4").Value = 0.54148292453891633
$ws.Range("I4").Value = 0.7390789588277753
